# LOQ4205.docx restructuring
#
# The document's content blocks got reshuffled (paragraph properties /
# pStyle stay put; only run text moves):
#   - the short PT/EN "Objetivos" paragraphs are replaced by the short
#     PT/EN "Programa resumido" lists
#   - the "Docente(s)" bullet gets the old PT objective text
#   - the old short PT/EN summary paragraphs become the long PT program
#     text / the EN objective text respectively
#   - the old long PT program paragraph becomes "Aulas Expositivas;
#     trabalhos e seminários."
#   - in the "Avaliação" bullet, the Método/Critério/Norma values shift
#     down one slot and the bibliography text becomes the new Norma
#     value
#   - the old bibliography paragraph becomes "5840535 - Messias Borges
#     Silva"
#
# Each single-run paragraph is rewritten by index via Paragraphs.Item(N)
# (stable regardless of execution order, since no paragraphs are added
# or removed). The one multi-run paragraph (Avaliação bullet) is edited
# by running Find scoped to that paragraph's own Range so only the
# targeted value run's text changes, leaving the bold labels and the
# other runs' formatting/line-breaks untouched.

$d = $word.ActiveDocument
$vbrk = [char]11

function Set-ParagraphText([int]$index, [string]$newText) {
    $d.Paragraphs.Item($index).Range.Text = $newText
}

function Replace-InParagraph([int]$index, [string]$oldText, [string]$newText) {
    $rng = $d.Paragraphs.Item($index).Range
    $found = $rng.Find.Execute($oldText, $true, $true, $false, $false, $false, `
                                $true, 1, $false, "", 0)
    if (-not $found) {
        throw "Find failed in paragraph $index for: $oldText"
    }
    $rng.Text = $newText
}

# 1) "Objetivos" (PT) paragraph -> short PT summary list.
$text1 = "1 – Planejamento e Controle da Qualidade" + $vbrk + `
    "2 – Melhoramentos da Produção" + $vbrk + `
    "3 – Desafios da produção" + $vbrk + `
    "4 – Controle da Qualidade"
Set-ParagraphText 6 $text1

# 2) "Objetivos" (EN, italic) paragraph -> short EN summary list.
$text2 = "1 - Quality Planning and Control" + $vbrk + `
    "2 - Production Improvements" + $vbrk + `
    "3 - Production challenges" + $vbrk + `
    "4 - Quality Control"
Set-ParagraphText 7 $text2

# 3) "Docente(s)" bullet -> old PT objective text.
$text3 = "Aprofundar os conceitos técnicos fundamentais de um curso de Engenharia de Produção, tendo em vista a sua formação generalista voltada para os mais diversos tipos de sistemas de produção."
Set-ParagraphText 9 $text3

# 4) Old short PT summary paragraph -> long PT program text.
$text4 = "1 – Planejamento e Controle da Qualidade" + $vbrk + `
    "Introdução. Planejamento e Controle da qualidade." + $vbrk + `
    "2 – Melhoramentos da Produção" + $vbrk + `
    "Introdução. Medidas e melhoramentos de desempenho. Prevenção e Recuperação de falhas. Administração da Qualidade Total." + $vbrk + `
    "3 – Desafios da produção" + $vbrk + `
    "Introdução. Tipo e formas de estratégias." + $vbrk + `
    "4 - CONTROLE DA QUALIDADE" + $vbrk + `
    "As Sete Ferramentas da Qualidade: Diagrama de Ishikawa, Histograma, Folha de Verificação, Estratificação, Diagrama de Pareto, Diagrama de Dispersão, Gráficos de Controle. Círculos de Controle da Qualidade"
Set-ParagraphText 11 $text4

# 5) Old short EN summary (italic) paragraph -> EN objective text.
$text5 = "To deepen the technical concepts of a Industrial Engineering course, in view of its general training aimed at the most diverse types of production systems."
Set-ParagraphText 12 $text5

# 6) Old PT long program paragraph -> "Aulas Expositivas; trabalhos e seminários."
$text6 = "Aulas Expositivas; trabalhos e seminários."
Set-ParagraphText 14 $text6

# 7) "Avaliação" bullet: Método/Critério/Norma values each shift down one
#    slot, and the bibliography text becomes the new Norma value.
#    Processed right-to-left (Norma, then Critério, then Método) so that
#    at each Find the searched-for text is still the single, unique
#    occurrence in the paragraph (earlier steps would otherwise
#    momentarily create a duplicate that a later Find could latch onto).
$oldNorma = "NF = (MF + PR)/2, onde PR é uma prova de recuperação."
$newNorma = "SLACK, N. et al. Administração da produção. São Paulo: Atlas, 2002. " + $vbrk + $vbrk + `
    "VENANZI, D; SILVA, O.R., Gerenciamento da Produçao e Operaçoes, LTC, 2014" + $vbrk + $vbrk + `
    "Textos complementares serão usados durante o curso."
Replace-InParagraph 17 $oldNorma $newNorma

$oldCriterio = "MF = (0,30*P1 + 0,30*P2 + 0,40*TRAB), onde P1 e P2 são provas e TRAB é a nota média de trabalhos e seminários."
$newCriterio = "NF = (MF + PR)/2, onde PR é uma prova de recuperação."
Replace-InParagraph 17 $oldCriterio $newCriterio

$oldMetodo = "Aulas Expositivas; trabalhos e seminários."
$newMetodo = "MF = (0,30*P1 + 0,30*P2 + 0,40*TRAB), onde P1 e P2 são provas e TRAB é a nota média de trabalhos e seminários."
Replace-InParagraph 17 $oldMetodo $newMetodo

# 8) Old bibliography paragraph -> "5840535 - Messias Borges Silva".
$text8 = "5840535 - Messias Borges Silva"
Set-ParagraphText 19 $text8
